$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new data point (class "apple") was measured and inserted right after the
# existing class-"apple" rows (as the new row 8), and a second new data point
# (class "orange") was inserted in the middle of the class-"orange" block (as
# the new row 12). Both insertions push the following rows down by one.
$ws.Rows.Item(8).Insert()
$ws.Rows.Item(12).Insert()

# Fill in the two newly inserted rows with their measured Hue/Saturation/Value
# and class label.
$ws.Range("A8").Value = 50.84785917930964
$ws.Range("B8").Value = 68.93794165563588
$ws.Range("C8").Value = 208.8774615619591
$ws.Range("D8").Value = "apple"

$ws.Range("A12").Value = 43.51049858109956
$ws.Range("B12").Value = 91.97278446389497
$ws.Range("C12").Value = 240.5210698167396
$ws.Range("D12").Value = "orange"

# Fix the typo in the "Class" label: "apples" -> "apple" (affects every row
# that carries that class, including the original rows and the new row 8).
$ws.Cells.Replace("apples", "apple")
